# Auto-generated cell updates derived from the canonical OOXML diff.
# Each (sheet,row) block updates the changed H:N profit/price columns in place.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1978
$ws.Range("I32").Value = 1973
$ws.Range("K32").Value = 1973
$ws.Range("M32").Value = -1647

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752

$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142

$ws.Range("H74").Value = 9899.333000000001
$ws.Range("I74").Value = 9899.333000000001
$ws.Range("K74").Value = 9899.333000000001
$ws.Range("M74").Value = -8963.333000000001

$ws.Range("H77").Value = 9899.333000000001
$ws.Range("I77").Value = 9899.333000000001
$ws.Range("K77").Value = 49496.665
$ws.Range("M77").Value = -44816.665

$ws.Range("H92").Value = 766.3182
$ws.Range("I92").Value = 862.0625
$ws.Range("K92").Value = 862.0625
$ws.Range("M92").Value = 385.9375

$ws.Range("H106").Value = 6250
$ws.Range("I106").Value = 6666.6665
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 6666.6665
$ws.Range("L106").Value = 5000
$ws.Range("M106").Value = -6035.6665
$ws.Range("N106").Value = -6262

$ws.Range("H138").Value = 3152.25
$ws.Range("J138").Value = 3849.44
$ws.Range("L138").Value = 11548.32
$ws.Range("N138").Value = -21828.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 17999
$ws.Range("J94").Value = 17999
$ws.Range("L94").Value = 17999
$ws.Range("N94").Value = -19801

$ws.Range("H132").Value = 1349.1714
$ws.Range("I132").Value = 1074.5
$ws.Range("K132").Value = 3223.5
$ws.Range("M132").Value = -693.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2063.6667
$ws.Range("I86").Value = 1829.5834
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1829.5834
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -706.5834
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 2063.6667
$ws.Range("I89").Value = 1829.5834
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 9147.916999999999
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -3531.916999999999
$ws.Range("N89").Value = -26232

$ws.Range("H134").Value = 2455.5881
$ws.Range("I134").Value = 2296.5625
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6889.6875
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4354.6875
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155

$ws.Range("H62").Value = 9998
$ws.Range("J62").Value = 9995
$ws.Range("L62").Value = 9995
$ws.Range("N62").Value = -11243

$ws.Range("H65").Value = 9998
$ws.Range("J65").Value = 9995
$ws.Range("L65").Value = 49975
$ws.Range("N65").Value = -56215

$ws.Range("H86").Value = 9177.6
$ws.Range("I86").Value = 9330.333000000001
$ws.Range("J86").Value = 8948.5
$ws.Range("K86").Value = 9330.333000000001
$ws.Range("L86").Value = 8948.5
$ws.Range("M86").Value = -8207.333000000001
$ws.Range("N86").Value = -11194.5

$ws.Range("H89").Value = 9177.6
$ws.Range("I89").Value = 9330.333000000001
$ws.Range("J89").Value = 8948.5
$ws.Range("K89").Value = 46651.665
$ws.Range("L89").Value = 44742.5
$ws.Range("M89").Value = -41035.665
$ws.Range("N89").Value = -55974.5

$ws.Range("H132").Value = 3874.8462
$ws.Range("I132").Value = 2964.8333
$ws.Range("J132").Value = 4654.857
$ws.Range("K132").Value = 8894.499899999999
$ws.Range("L132").Value = 13964.571
$ws.Range("M132").Value = -6364.499899999999
$ws.Range("N132").Value = -19024.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 831.5833
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 831.5833
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2494.7499
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2718.7499

$ws.Range("H34").Value = 6581.5835
$ws.Range("J34").Value = 6997.778
$ws.Range("L34").Value = 20993.334
$ws.Range("N34").Value = -21161.334

$ws.Range("H39").Value = 14854.286
$ws.Range("J39").Value = 14854.286
$ws.Range("L39").Value = 44562.858
$ws.Range("N39").Value = -45150.858

$ws.Range("H55").Value = 9941.666999999999
$ws.Range("J55").Value = 13662.5
$ws.Range("L55").Value = 40987.5
$ws.Range("N55").Value = -41341.5

$ws.Range("H113").Value = 2186.4285
$ws.Range("J113").Value = 2186.4285
$ws.Range("L113").Value = 6559.2855
$ws.Range("N113").Value = -10899.2855

$ws.Range("H132").Value = 1405.1
$ws.Range("I132").Value = 1150.2858
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 10352.5722
$ws.Range("L132").Value = 17996.9994
$ws.Range("M132").Value = -7822.572200000001
$ws.Range("N132").Value = -23056.9994

$ws.Range("H135").Value = 831.5833
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 831.5833
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 7484.2497
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -12554.2497

$ws.Range("H139").Value = 1414.5
$ws.Range("I139").Value = 1414.5
$ws.Range("K139").Value = 4243.5
$ws.Range("M139").Value = 896.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 14250
$ws.Range("J28").Value = 14250
$ws.Range("L28").Value = 14250
$ws.Range("N28").Value = -14634

$ws.Range("H80").Value = 3564.75
$ws.Range("J80").Value = 3631.5
$ws.Range("L80").Value = 3631.5
$ws.Range("N80").Value = -5627.5

$ws.Range("H83").Value = 3564.75
$ws.Range("J83").Value = 3631.5
$ws.Range("L83").Value = 18157.5
$ws.Range("N83").Value = -28141.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2755
$ws.Range("I7").Value = 2755
$ws.Range("K7").Value = 2755
$ws.Range("M7").Value = -2643

$ws.Range("H100").Value = 2816.8333
$ws.Range("I100").Value = 2380.2
$ws.Range("K100").Value = 2380.2
$ws.Range("M100").Value = -1839.2

$ws.Range("H126").Value = 2755
$ws.Range("I126").Value = 2755
$ws.Range("K126").Value = 8265
$ws.Range("M126").Value = -5795

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1446.8182
$ws.Range("I126").Value = 771.6667
$ws.Range("K126").Value = 2315.0001
$ws.Range("M126").Value = 154.9998999999998

$ws.Range("H132").Value = 3240
$ws.Range("I132").Value = 3014.3333
$ws.Range("J132").Value = 3465.6667
$ws.Range("K132").Value = 9042.999899999999
$ws.Range("L132").Value = 10397.0001
$ws.Range("M132").Value = -6512.999899999999
$ws.Range("N132").Value = -15457.0001

Write-Output "Updated 174 cells, cleared 3 cells."
